$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 128: remove D128 ---
$ws.Range("D128").Value = $null

# --- Rows 129-132: update D values (reduce float precision noise) ---
$ws.Range("D129").Value = 0.293130374
$ws.Range("D130").Value = 0.30455641
$ws.Range("D131").Value = 0.224211961
$ws.Range("D132").Value = 0.2859344540000001

# --- Rows 133-139: update C values (reduce float precision noise) ---
$ws.Range("C133").Value = 0.107440335
$ws.Range("C134").Value = -0.634467653
$ws.Range("C135").Value = -0.255775563
$ws.Range("C136").Value = 0.03129978100000003
$ws.Range("C137").Value = 0.05894410600000002
$ws.Range("C138").Value = -0.18190582
$ws.Range("C139").Value = -0.09604011399999998

# --- Rows 136-139: add new D values ---
$ws.Range("D136").Value = 0.868324474
$ws.Range("D137").Value = 0.864449887
$ws.Range("D138").Value = 0.6740699019999999
$ws.Range("D139").Value = 0.696856487

# --- Row 140: update B value, add C value ---
$ws.Range("B140").Value = -0.387365295
$ws.Range("C140").Value = 0.318216844

# --- New rows 141-145 (ifoCAST sampling extension) ---
# Copy formatting from the last existing label cell (A140) so the new
# date-label cells pick up the same bold/border/centered style (s="1")
# without introducing a brand-new style entry in styles.xml.

$ws.Range("A141").Value = "2025-07-25_diff"
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A141").PasteSpecial(-4122) | Out-Null
$ws.Range("B141").Value = -0.431278794
$ws.Range("C141").Value = 0.407980578

$ws.Range("A142").Value = "2025-08-07_diff"
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A142").PasteSpecial(-4122) | Out-Null
$ws.Range("C142").Value = 0.230511215

$ws.Range("A143").Value = "2025-08-22_diff"
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A143").PasteSpecial(-4122) | Out-Null
$ws.Range("C143").Value = 0.243228468

$ws.Range("A144").Value = "2025-08-25_diff"
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A144").PasteSpecial(-4122) | Out-Null
$ws.Range("C144").Value = 0.118688665

$ws.Range("A145").Value = "2025-09-08_diff"
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A145").PasteSpecial(-4122) | Out-Null
$ws.Range("C145").Value = 0.343027307

$excel.CutCopyMode = $false
